$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "system, System, backup@backdoor.com" = "System, backup@backdoor.com, system"
    "System, dnasr281@gmail.com" = "dnasr281@gmail.com, System"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
}

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value2 = $map[$val]
    }
}
